# B1/B2 PowerPoint update
#
#  1) Re-style the financial-documents comparison table (slide 5) from the
#     deck's bespoke "Table_0" style to the built-in
#     "Medium Style 2 - Accent 1" table style.
#  2) Switch the deck's theme colour palette from the custom
#     "Integral / Red Violet" scheme over to the standard "Office" scheme
#     (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).

$p = $ppt.ActivePresentation

# --- 1) Table style --------------------------------------------------------
$slide      = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$table      = $tableShape.Table
$table.ApplyStyle("{9DF418A3-52CC-4C8A-9F19-C31BA865F6CE}")

# --- 2) Theme colours -------------------------------------------------------
# RGB() values below are the standard VBA long (R + G*256 + B*65536) for the
# Office theme's twelve scheme colours.
$colorScheme = $p.Slides.Item(1).ThemeColorScheme
$colorScheme.Item(1).RGB  = 0          # dk1      000000
$colorScheme.Item(2).RGB  = 16777215   # lt1      FFFFFF
$colorScheme.Item(3).RGB  = 6968388    # dk2      44546A
$colorScheme.Item(4).RGB  = 15132391   # lt2      E7E6E6
$colorScheme.Item(5).RGB  = 13998939   # accent1  5B9BD5
$colorScheme.Item(6).RGB  = 3243501    # accent2  ED7D31
$colorScheme.Item(7).RGB  = 10855845   # accent3  A5A5A5
$colorScheme.Item(8).RGB  = 49407      # accent4  FFC000
$colorScheme.Item(9).RGB  = 12874308   # accent5  4472C4
$colorScheme.Item(10).RGB = 4697456    # accent6  70AD47
$colorScheme.Item(11).RGB = 12673797   # hlink    0563C1
$colorScheme.Item(12).RGB = 7491477    # folHlink 954F72
